# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for each coin row to reflect the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.169.66"
$ws.Range("E2").Value = "  +3.26%  "
$ws.Range("D3").Value = "3.116.21"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'523.23"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").Value = "'145.03"
$ws.Range("E6").Value = "  +2.84%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.441"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").Value = "'7.39"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("E11").Value = "  +3.69%  "
$ws.Range("D12").Value = "3.649.14"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").Value = "'27.30"
$ws.Range("E14").Value = "  +7.34%  "
$ws.Range("D15").Value = "'0.0000167"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "59.139.75"
$ws.Range("E16").Value = "  +3.17%  "
$ws.Range("D17").Value = "3.118.78"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("E18").Value = "  +4.79%  "
$ws.Range("D19").Value = "'13.11"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "'8.29"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("D21").Value = "'345.36"
$ws.Range("E21").Value = "  +2.46%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  +2.80%  "
$ws.Range("D24").Value = "'65.87"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "0.0₃0939"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'6.77"
$ws.Range("E28").Value = "  +5.67%  "
$ws.Range("D29").Value = "'7.34"
$ws.Range("E29").Value = "  +3.96%  "
$ws.Range("D30").Value = "'1.85"
$ws.Range("E30").Value = "  +2.85%  "
$ws.Range("E31").Value = "  +4.49%  "
$ws.Range("D32").Value = "'21.17"
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("D33").Value = "'155.66"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "'4.69"
$ws.Range("E34").Value = "  +3.76%  "
$ws.Range("E35").Value = "  +6.21%  "
$ws.Range("D36").Value = "'27.39"
$ws.Range("E36").Value = "  +5.81%  "
$ws.Range("E37").Value = "  +6.68%  "
$ws.Range("D38").Value = "'0.0690"
$ws.Range("E38").Value = "  +2.87%  "
$ws.Range("D39").Value = "'3.97"
$ws.Range("E39").Value = "  +3.33%  "
$ws.Range("D40").Value = "3.157.58"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").Value = "'36.96"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'0.667"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "'1.47"
$ws.Range("E44").Value = "  +6.47%  "
$ws.Range("D45").Value = "2.289.33"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("E46").Value = "  +2.97%  "
$ws.Range("D47").Value = "'21.25"
$ws.Range("E47").Value = "  +6.22%  "
$ws.Range("D48").Value = "'0.968"
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("D49").Value = "'6.05"
$ws.Range("E49").Value = "  +3.78%  "
$ws.Range("D50").Value = "'0.762"
$ws.Range("E50").Value = "  +12.04%  "
$ws.Range("D51").Value = "'263.84"
$ws.Range("E51").Value = "  +11.74%  "
